# Applies the cryptos-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.887.19'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').Value = '1.830.63'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.56'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6942'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07667'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3043'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.26'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07807'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '92.84'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.834.26'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.090'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6826'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008235'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').Value = '28.910.91'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('D20').Value = '2.074.28'
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.66'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.449'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.729'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.539'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.224'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.139'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05116'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7723'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.848'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.140'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.693'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = '1.273.36'
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01857'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.700'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9522'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.131'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '106.87'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9991'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.693'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5165'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = '1.973.81'
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '63.67'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -7.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.750'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.966'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.67%  '
